$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "speedup" header merged over E10:F10, centered (same style as H1:L1)
$ws.Range("E10").Value = "speedup"
$ws.Range("E10:F10").Merge() | Out-Null
$ws.Range("E10:F10").HorizontalAlignment = -4108  # xlCenter

# Row 11: column-size headers
$ws.Range("B11").Value = 1024
$ws.Range("C11").Value = 4096
$ws.Range("E11").Value = 1024
$ws.Range("F11").Value = 4096

# Row 12: serial baseline timings
$ws.Range("A12").Value = "serial"
$ws.Range("B12").Value = 0.33892099999999997
$ws.Range("C12").Value = 23.258996

# Row 13: 2 threads
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 0.16855100000000001
$ws.Range("C13").Value = 12.561260000000001
$ws.Range("E13").Formula = "=B12/B13"
$ws.Range("F13").Formula = "=C12/C13"
$ws.Range("H13").Formula = "=AVERAGE(E13:F13)"

# Row 14: 5 threads
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = 0.082139
$ws.Range("C14").Value = 5.782565
$ws.Range("E14").Formula = "=B12/B14"
$ws.Range("F14").Formula = "=C12/C14"

# Row 15: 8 threads
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 0.065322
$ws.Range("C15").Value = 4.4164570000000003
$ws.Range("E15").Formula = "=B12/B15"
$ws.Range("F15").Formula = "=C12/C15"

# Row 16: 10 threads
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = 0.069593
$ws.Range("C16").Value = 4.2484320000000002
$ws.Range("E16").Formula = "=B12/B16"
$ws.Range("F16").Formula = "=C12/C16"

# Row 17: 12 threads
$ws.Range("A17").Value = 12
$ws.Range("B17").Value = 0.068674
$ws.Range("C17").Value = 4.2780319999999996
$ws.Range("E17").Formula = "=B12/B17"
$ws.Range("F17").Formula = "=C12/C17"

# H14:H17 filled together so Excel records it as one shared formula group
$ws.Range("H14:H17").Formula = "=AVERAGE(E14:F14)"

$ws.Range("H13").Select() | Out-Null
